$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at D (shifts old D -> G) to make room for
# Corequisites (D), Concurrent (E), Recommended (F); old D becomes G.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Cells.Item(1, 4).Value = 'Corequisites'
$ws.Cells.Item(1, 5).Value = 'Concurrent'
$ws.Cells.Item(1, 6).Value = 'Recommended'
# G1 ("Terms Typically Offered") already shifted correctly from old D1.

# Per-row data
$ws.Cells.Item(2, 4).Value = 'NA'
$ws.Cells.Item(2, 5).Value = 'NA'
$ws.Cells.Item(2, 6).Value = 'NA'
$ws.Cells.Item(2, 7).Value = 'F'
$ws.Cells.Item(3, 3).Value = 'PSY 201 or PSY 202 and consent of department head.'
$ws.Cells.Item(3, 4).Value = 'NA'
$ws.Cells.Item(3, 5).Value = 'NA'
$ws.Cells.Item(3, 6).Value = 'NA'
$ws.Cells.Item(3, 7).Value = 'TBD'
$ws.Cells.Item(4, 4).Value = 'NA'
$ws.Cells.Item(4, 5).Value = 'NA'
$ws.Cells.Item(4, 6).Value = 'NA'
$ws.Cells.Item(4, 7).Value = 'F, W, SP'
$ws.Cells.Item(5, 4).Value = 'NA'
$ws.Cells.Item(5, 5).Value = 'NA'
$ws.Cells.Item(5, 6).Value = 'NA'
$ws.Cells.Item(5, 7).Value = 'F, W'
$ws.Cells.Item(6, 4).Value = 'NA'
$ws.Cells.Item(6, 5).Value = 'NA'
$ws.Cells.Item(6, 6).Value = 'NA'
$ws.Cells.Item(6, 7).Value = 'TBD'
$ws.Cells.Item(7, 4).Value = 'NA'
$ws.Cells.Item(7, 5).Value = 'NA'
$ws.Cells.Item(7, 6).Value = 'NA'
$ws.Cells.Item(7, 7).Value = 'TBD'
$ws.Cells.Item(8, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(8, 4).Value = 'NA'
$ws.Cells.Item(8, 5).Value = 'NA'
$ws.Cells.Item(8, 6).Value = 'NA'
$ws.Cells.Item(8, 7).Value = 'F, W'
$ws.Cells.Item(9, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(9, 4).Value = 'NA'
$ws.Cells.Item(9, 5).Value = 'NA'
$ws.Cells.Item(9, 6).Value = 'NA'
$ws.Cells.Item(9, 7).Value = 'F, W, SP'
$ws.Cells.Item(10, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(10, 4).Value = 'NA'
$ws.Cells.Item(10, 5).Value = 'NA'
$ws.Cells.Item(10, 6).Value = 'NA'
$ws.Cells.Item(10, 7).Value = 'F, W, SP'
$ws.Cells.Item(11, 4).Value = 'NA'
$ws.Cells.Item(11, 5).Value = 'NA'
$ws.Cells.Item(11, 6).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(11, 7).Value = 'TBD'
$ws.Cells.Item(12, 4).Value = 'NA'
$ws.Cells.Item(12, 5).Value = 'NA'
$ws.Cells.Item(12, 6).Value = 'NA'
$ws.Cells.Item(12, 7).Value = 'TBD'
$ws.Cells.Item(13, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(13, 4).Value = 'NA'
$ws.Cells.Item(13, 5).Value = 'NA'
$ws.Cells.Item(13, 6).Value = 'NA'
$ws.Cells.Item(13, 7).Value = 'TBD'
$ws.Cells.Item(14, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(14, 4).Value = 'NA'
$ws.Cells.Item(14, 5).Value = 'NA'
$ws.Cells.Item(14, 6).Value = 'NA'
$ws.Cells.Item(14, 7).Value = 'F, W'
$ws.Cells.Item(15, 3).Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D1, D3, or D4/E.'
$ws.Cells.Item(15, 4).Value = 'NA'
$ws.Cells.Item(15, 5).Value = 'NA'
$ws.Cells.Item(15, 6).Value = 'Completion of USCP requirement.'
$ws.Cells.Item(15, 7).Value = 'SP '
$ws.Cells.Item(16, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(16, 4).Value = 'NA'
$ws.Cells.Item(16, 5).Value = 'NA'
$ws.Cells.Item(16, 6).Value = 'NA'
$ws.Cells.Item(16, 7).Value = 'F, W, SP'
$ws.Cells.Item(17, 3).Value = 'CD/PSY 256 or CD/EDUC 207.'
$ws.Cells.Item(17, 4).Value = 'NA'
$ws.Cells.Item(17, 5).Value = 'NA'
$ws.Cells.Item(17, 6).Value = 'NA'
$ws.Cells.Item(17, 7).Value = 'F, W, SP'
$ws.Cells.Item(18, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(18, 4).Value = 'NA'
$ws.Cells.Item(18, 5).Value = 'NA'
$ws.Cells.Item(18, 6).Value = 'NA'
$ws.Cells.Item(18, 7).Value = 'TBD'
$ws.Cells.Item(19, 3).Value = 'Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and completion of GE Area D4/E.'
$ws.Cells.Item(19, 4).Value = 'NA'
$ws.Cells.Item(19, 5).Value = 'NA'
$ws.Cells.Item(19, 6).Value = 'PSY 201 or PSY 202 (GE Area D4/E).'
$ws.Cells.Item(19, 7).Value = 'TBD '
$ws.Cells.Item(20, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(20, 4).Value = 'NA'
$ws.Cells.Item(20, 5).Value = 'NA'
$ws.Cells.Item(20, 6).Value = 'NA'
$ws.Cells.Item(20, 7).Value = 'TBD'
$ws.Cells.Item(21, 3).Value = 'Junior standing or Psychology major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Cells.Item(21, 4).Value = 'NA'
$ws.Cells.Item(21, 5).Value = 'NA'
$ws.Cells.Item(21, 6).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(21, 7).Value = 'F, W, SP '
$ws.Cells.Item(22, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(22, 4).Value = 'NA'
$ws.Cells.Item(22, 5).Value = 'NA'
$ws.Cells.Item(22, 6).Value = 'NA'
$ws.Cells.Item(22, 7).Value = 'W'
$ws.Cells.Item(23, 4).Value = 'NA'
$ws.Cells.Item(23, 5).Value = 'NA'
$ws.Cells.Item(23, 6).Value = 'NA'
$ws.Cells.Item(23, 7).Value = 'F, W, SP'
$ws.Cells.Item(24, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(24, 4).Value = 'NA'
$ws.Cells.Item(24, 5).Value = 'NA'
$ws.Cells.Item(24, 6).Value = 'NA'
$ws.Cells.Item(24, 7).Value = 'W, SP'
$ws.Cells.Item(25, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(25, 4).Value = 'NA'
$ws.Cells.Item(25, 5).Value = 'NA'
$ws.Cells.Item(25, 6).Value = 'NA'
$ws.Cells.Item(25, 7).Value = 'W'
$ws.Cells.Item(26, 3).Value = 'PSY 201 or PSY 202; and STAT 217.'
$ws.Cells.Item(26, 4).Value = 'NA'
$ws.Cells.Item(26, 5).Value = 'NA'
$ws.Cells.Item(26, 6).Value = 'NA'
$ws.Cells.Item(26, 7).Value = 'F, W, SP'
$ws.Cells.Item(27, 4).Value = 'NA'
$ws.Cells.Item(27, 5).Value = 'NA'
$ws.Cells.Item(27, 6).Value = 'NA'
$ws.Cells.Item(27, 7).Value = 'F, W'
$ws.Cells.Item(28, 3).Value = 'PSY 329; and STAT 217.'
$ws.Cells.Item(28, 4).Value = 'NA'
$ws.Cells.Item(28, 5).Value = 'NA'
$ws.Cells.Item(28, 6).Value = 'NA'
$ws.Cells.Item(28, 7).Value = 'F, W, SP'
$ws.Cells.Item(29, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(29, 4).Value = 'NA'
$ws.Cells.Item(29, 5).Value = 'NA'
$ws.Cells.Item(29, 6).Value = 'NA'
$ws.Cells.Item(29, 7).Value = 'F, W, SP'
$ws.Cells.Item(30, 3).Value = 'PSY 201 or PSY 202; and one of the ASCI 112, BIO 111, BIO 123, BIO 161, or BIO 213.'
$ws.Cells.Item(30, 4).Value = 'NA'
$ws.Cells.Item(30, 5).Value = 'NA'
$ws.Cells.Item(30, 6).Value = 'STAT 217 or STAT 218.'
$ws.Cells.Item(30, 7).Value = 'F, W '
$ws.Cells.Item(31, 4).Value = 'NA'
$ws.Cells.Item(31, 5).Value = 'NA'
$ws.Cells.Item(31, 6).Value = 'NA'
$ws.Cells.Item(31, 7).Value = 'F, W, SP'
$ws.Cells.Item(32, 3).Value = 'Junior standing or Psychology major; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one lower-division course in GE Area D.'
$ws.Cells.Item(32, 4).Value = 'NA'
$ws.Cells.Item(32, 5).Value = 'NA'
$ws.Cells.Item(32, 6).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(32, 7).Value = 'SP '
$ws.Cells.Item(33, 3).Value = 'PSY 201 or PSY 202; CD/PSY 256 or CD 305; and junior standing.'
$ws.Cells.Item(33, 4).Value = 'NA'
$ws.Cells.Item(33, 5).Value = 'NA'
$ws.Cells.Item(33, 6).Value = 'NA'
$ws.Cells.Item(33, 7).Value = 'W, SP'
$ws.Cells.Item(34, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(34, 4).Value = 'NA'
$ws.Cells.Item(34, 5).Value = 'NA'
$ws.Cells.Item(34, 6).Value = 'NA'
$ws.Cells.Item(34, 7).Value = 'W'
$ws.Cells.Item(35, 3).Value = 'PSY 252.'
$ws.Cells.Item(35, 4).Value = 'NA'
$ws.Cells.Item(35, 5).Value = 'NA'
$ws.Cells.Item(35, 6).Value = 'NA'
$ws.Cells.Item(35, 7).Value = 'SP'
$ws.Cells.Item(36, 3).Value = 'PSY 201 or PSY 202 and at least one other PSY course.'
$ws.Cells.Item(36, 4).Value = 'NA'
$ws.Cells.Item(36, 5).Value = 'NA'
$ws.Cells.Item(36, 6).Value = 'NA'
$ws.Cells.Item(36, 7).Value = 'W, SP'
$ws.Cells.Item(37, 3).Value = 'PSY 201 or PSY 202; and sophomore standing.'
$ws.Cells.Item(37, 4).Value = 'NA'
$ws.Cells.Item(37, 5).Value = 'NA'
$ws.Cells.Item(37, 6).Value = 'NA'
$ws.Cells.Item(37, 7).Value = 'F, W, SP'
$ws.Cells.Item(38, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(38, 4).Value = 'NA'
$ws.Cells.Item(38, 5).Value = 'NA'
$ws.Cells.Item(38, 6).Value = 'NA'
$ws.Cells.Item(38, 7).Value = 'F, W, SP'
$ws.Cells.Item(39, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(39, 4).Value = 'NA'
$ws.Cells.Item(39, 5).Value = 'NA'
$ws.Cells.Item(39, 6).Value = 'NA'
$ws.Cells.Item(39, 7).Value = 'TBD'
$ws.Cells.Item(40, 4).Value = 'NA'
$ws.Cells.Item(40, 5).Value = 'NA'
$ws.Cells.Item(40, 6).Value = 'NA'
$ws.Cells.Item(40, 7).Value = 'TBD'
$ws.Cells.Item(41, 3).Value = 'PSY 201 or PSY 202.'
$ws.Cells.Item(41, 4).Value = 'NA'
$ws.Cells.Item(41, 5).Value = 'NA'
$ws.Cells.Item(41, 6).Value = 'NA'
$ws.Cells.Item(41, 7).Value = 'F, W, SP'
$ws.Cells.Item(42, 3).Value = 'PSY 333.'
$ws.Cells.Item(42, 4).Value = 'NA'
$ws.Cells.Item(42, 5).Value = 'NA'
$ws.Cells.Item(42, 6).Value = 'NA'
$ws.Cells.Item(42, 7).Value = 'TBD'
$ws.Cells.Item(43, 3).Value = 'CD 304, CD 305 or CD 306; or CD/PSY 256 and PSY 305.'
$ws.Cells.Item(43, 4).Value = 'NA'
$ws.Cells.Item(43, 5).Value = 'NA'
$ws.Cells.Item(43, 6).Value = 'NA'
$ws.Cells.Item(43, 7).Value = 'F, W, SP'
$ws.Cells.Item(44, 3).Value = 'PSY 256 or PSY 305 and senior standing.'
$ws.Cells.Item(44, 4).Value = 'NA'
$ws.Cells.Item(44, 5).Value = 'NA'
$ws.Cells.Item(44, 6).Value = 'NA'
$ws.Cells.Item(44, 7).Value = 'TBD'
$ws.Cells.Item(45, 3).Value = 'PSY 340.'
$ws.Cells.Item(45, 4).Value = 'NA'
$ws.Cells.Item(45, 5).Value = 'NA'
$ws.Cells.Item(45, 6).Value = 'NA'
$ws.Cells.Item(45, 7).Value = 'W, SP'
$ws.Cells.Item(46, 3).Value = 'CD 304 or CD 305 or PSY 419 or PSY 420 or PSY 421; and CD 329 or PSY 329.'
$ws.Cells.Item(46, 4).Value = 'NA'
$ws.Cells.Item(46, 5).Value = 'NA'
$ws.Cells.Item(46, 6).Value = 'NA'
$ws.Cells.Item(46, 7).Value = 'F, W'
$ws.Cells.Item(47, 3).Value = 'PSY 333.'
$ws.Cells.Item(47, 4).Value = 'NA'
$ws.Cells.Item(47, 5).Value = 'NA'
$ws.Cells.Item(47, 6).Value = 'NA'
$ws.Cells.Item(47, 7).Value = 'TBD'
$ws.Cells.Item(48, 3).Value = 'PSY 201 or PSY 202; and PSY 329.'
$ws.Cells.Item(48, 4).Value = 'NA'
$ws.Cells.Item(48, 5).Value = 'NA'
$ws.Cells.Item(48, 6).Value = 'NA'
$ws.Cells.Item(48, 7).Value = 'F, W, SP'
$ws.Cells.Item(49, 3).Value = 'PSY 329, PSY 333, Psychology and Child Development majors only, junior standing, and consent of instructor.'
$ws.Cells.Item(49, 4).Value = 'NA'
$ws.Cells.Item(49, 5).Value = 'NA'
$ws.Cells.Item(49, 6).Value = 'NA'
$ws.Cells.Item(49, 7).Value = 'F, W, SP'
$ws.Cells.Item(50, 3).Value = 'PSY 329, PSY 333, Psychology and Child Development majors only, junior standing, and consent of instructor.'
$ws.Cells.Item(50, 4).Value = 'NA'
$ws.Cells.Item(50, 5).Value = 'NA'
$ws.Cells.Item(50, 6).Value = 'NA'
$ws.Cells.Item(50, 7).Value = 'F, W, SP'
$ws.Cells.Item(51, 3).Value = 'PSY 254, or graduate standing.'
$ws.Cells.Item(51, 4).Value = 'NA'
$ws.Cells.Item(51, 5).Value = 'NA'
$ws.Cells.Item(51, 6).Value = 'NA'
$ws.Cells.Item(51, 7).Value = 'TBD'
$ws.Cells.Item(52, 3).Value = 'PSY 323, Psychology and Child Development majors only, junior standing, and consent of instructor.'
$ws.Cells.Item(52, 4).Value = 'NA'
$ws.Cells.Item(52, 5).Value = 'NA'
$ws.Cells.Item(52, 6).Value = 'NA'
$ws.Cells.Item(52, 7).Value = 'F, W, SP'
$ws.Cells.Item(53, 3).Value = 'PSY 323, Psychology and Child Development majors only, junior standing, and consent of instructor.'
$ws.Cells.Item(53, 4).Value = 'NA'
$ws.Cells.Item(53, 5).Value = 'NA'
$ws.Cells.Item(53, 6).Value = 'NA'
$ws.Cells.Item(53, 7).Value = 'F, W, SP'
$ws.Cells.Item(54, 3).Value = 'PSY 333.'
$ws.Cells.Item(54, 4).Value = 'NA'
$ws.Cells.Item(54, 5).Value = 'NA'
$ws.Cells.Item(54, 6).Value = 'NA'
$ws.Cells.Item(54, 7).Value = 'F, W'
$ws.Cells.Item(55, 3).Value = 'PSY 201 or PSY 202 and junior standing.'
$ws.Cells.Item(55, 4).Value = 'NA'
$ws.Cells.Item(55, 5).Value = 'NA'
$ws.Cells.Item(55, 6).Value = 'NA'
$ws.Cells.Item(55, 7).Value = 'F, W, SP'
$ws.Cells.Item(56, 3).Value = 'Senior standing; PSY 329; Psychology and Child Development majors only.'
$ws.Cells.Item(56, 4).Value = 'NA'
$ws.Cells.Item(56, 5).Value = 'NA'
$ws.Cells.Item(56, 6).Value = 'NA'
$ws.Cells.Item(56, 7).Value = 'F, W, SP'
$ws.Cells.Item(57, 3).Value = 'PSY 461; Psychology and Child Development majors only.'
$ws.Cells.Item(57, 4).Value = 'NA'
$ws.Cells.Item(57, 5).Value = 'NA'
$ws.Cells.Item(57, 6).Value = 'NA'
$ws.Cells.Item(57, 7).Value = 'F, W, SP'
$ws.Cells.Item(58, 3).Value = 'PSY 201 or PSY 202 and junior standing.'
$ws.Cells.Item(58, 4).Value = 'NA'
$ws.Cells.Item(58, 5).Value = 'NA'
$ws.Cells.Item(58, 6).Value = 'NA'
$ws.Cells.Item(58, 7).Value = 'F, SP'
$ws.Cells.Item(59, 4).Value = 'NA'
$ws.Cells.Item(59, 5).Value = 'NA'
$ws.Cells.Item(59, 6).Value = 'NA'
$ws.Cells.Item(59, 7).Value = 'TBD'
$ws.Cells.Item(60, 3).Value = 'PSY 252 or PSY 254 or PSY 256.'
$ws.Cells.Item(60, 4).Value = 'NA'
$ws.Cells.Item(60, 5).Value = 'NA'
$ws.Cells.Item(60, 6).Value = 'NA'
$ws.Cells.Item(60, 7).Value = 'SP'
$ws.Cells.Item(61, 3).Value = 'PSY 340.'
$ws.Cells.Item(61, 4).Value = 'NA'
$ws.Cells.Item(61, 5).Value = 'NA'
$ws.Cells.Item(61, 6).Value = 'NA'
$ws.Cells.Item(61, 7).Value = 'F, W'
$ws.Cells.Item(62, 4).Value = 'NA'
$ws.Cells.Item(62, 5).Value = 'NA'
$ws.Cells.Item(62, 6).Value = 'NA'
$ws.Cells.Item(62, 7).Value = 'TBD'
$ws.Cells.Item(63, 4).Value = 'NA'
$ws.Cells.Item(63, 5).Value = 'NA'
$ws.Cells.Item(63, 6).Value = 'NA'
$ws.Cells.Item(63, 7).Value = 'TBD'
$ws.Cells.Item(64, 4).Value = 'NA'
$ws.Cells.Item(64, 5).Value = 'NA'
$ws.Cells.Item(64, 6).Value = 'NA'
$ws.Cells.Item(64, 7).Value = 'F, W, SP'
$ws.Cells.Item(65, 4).Value = 'NA'
$ws.Cells.Item(65, 5).Value = 'NA'
$ws.Cells.Item(65, 6).Value = 'NA'
$ws.Cells.Item(65, 7).Value = 'SP'
$ws.Cells.Item(66, 4).Value = 'NA'
$ws.Cells.Item(66, 5).Value = 'NA'
$ws.Cells.Item(66, 6).Value = 'NA'
$ws.Cells.Item(66, 7).Value = 'F'
$ws.Cells.Item(67, 4).Value = 'NA'
$ws.Cells.Item(67, 5).Value = 'NA'
$ws.Cells.Item(67, 6).Value = 'NA'
$ws.Cells.Item(67, 7).Value = 'F'
$ws.Cells.Item(68, 3).Value = 'PSY 560 and admission to MS Psychology program.'
$ws.Cells.Item(68, 4).Value = 'NA'
$ws.Cells.Item(68, 5).Value = 'NA'
$ws.Cells.Item(68, 6).Value = 'NA'
$ws.Cells.Item(68, 7).Value = 'W'
$ws.Cells.Item(69, 3).Value = 'PSY 555, PSY 560 and admission to MS Psychology program.'
$ws.Cells.Item(69, 4).Value = 'NA'
$ws.Cells.Item(69, 5).Value = 'NA'
$ws.Cells.Item(69, 6).Value = 'NA'
$ws.Cells.Item(69, 7).Value = 'SP'
$ws.Cells.Item(70, 4).Value = 'NA'
$ws.Cells.Item(70, 5).Value = 'NA'
$ws.Cells.Item(70, 6).Value = 'NA'
$ws.Cells.Item(70, 7).Value = 'F'
$ws.Cells.Item(71, 3).Value = 'PSY 520; PSY 560; and admission to MS Psychology program.'
$ws.Cells.Item(71, 4).Value = 'NA'
$ws.Cells.Item(71, 5).Value = 'NA'
$ws.Cells.Item(71, 6).Value = 'NA'
$ws.Cells.Item(71, 7).Value = 'SU'
$ws.Cells.Item(72, 3).Value = 'PSY 560 and admission to MS Psychology program.'
$ws.Cells.Item(72, 4).Value = 'NA'
$ws.Cells.Item(72, 5).Value = 'NA'
$ws.Cells.Item(72, 6).Value = 'NA'
$ws.Cells.Item(72, 7).Value = 'W'
$ws.Cells.Item(73, 3).Value = 'PSY 560 and admission to MS Psychology program.'
$ws.Cells.Item(73, 4).Value = 'NA'
$ws.Cells.Item(73, 5).Value = 'NA'
$ws.Cells.Item(73, 6).Value = 'NA'
$ws.Cells.Item(73, 7).Value = 'SP'
$ws.Cells.Item(74, 3).Value = 'PSY 555, PSY 560, PSY 565 and admission to MS Psychology program, or consent of instructor.'
$ws.Cells.Item(74, 4).Value = 'NA'
$ws.Cells.Item(74, 5).Value = 'NA'
$ws.Cells.Item(74, 6).Value = 'NA'
$ws.Cells.Item(74, 7).Value = 'F'
$ws.Cells.Item(75, 3).Value = 'PSY 520; PSY 555; PSY 560; PSY 565; and admission to MS Psychology program.'
$ws.Cells.Item(75, 4).Value = 'NA'
$ws.Cells.Item(75, 5).Value = 'NA'
$ws.Cells.Item(75, 6).Value = 'NA'
$ws.Cells.Item(75, 7).Value = 'F, W, SP'
$ws.Cells.Item(76, 4).Value = 'NA'
$ws.Cells.Item(76, 5).Value = 'NA'
$ws.Cells.Item(76, 6).Value = 'NA'
$ws.Cells.Item(76, 7).Value = 'TBD'
$ws.Cells.Item(77, 3).Value = 'PSY 520; PSY 555; and admission to MS Psychology program.'
$ws.Cells.Item(77, 4).Value = 'NA'
$ws.Cells.Item(77, 5).Value = 'NA'
$ws.Cells.Item(77, 6).Value = 'NA'
$ws.Cells.Item(77, 7).Value = 'F'
$ws.Cells.Item(78, 3).Value = 'PSY 535; PSY 555; PSY 560; and admission to MS Psychology program.'
$ws.Cells.Item(78, 4).Value = 'NA'
$ws.Cells.Item(78, 5).Value = 'NA'
$ws.Cells.Item(78, 6).Value = 'NA'
$ws.Cells.Item(78, 7).Value = 'W'
$ws.Cells.Item(79, 4).Value = 'NA'
$ws.Cells.Item(79, 5).Value = 'NA'
$ws.Cells.Item(79, 6).Value = 'NA'
$ws.Cells.Item(79, 7).Value = 'W'
$ws.Cells.Item(80, 3).Value = 'PSY 520; PSY 560; and admission to MS Psychology program.'
$ws.Cells.Item(80, 4).Value = 'NA'
$ws.Cells.Item(80, 5).Value = 'NA'
$ws.Cells.Item(80, 6).Value = 'NA'
$ws.Cells.Item(80, 7).Value = 'SP'
$ws.Cells.Item(81, 3).Value = 'PSY 569, PSY 564 and consent of MS program committee.'
$ws.Cells.Item(81, 4).Value = 'NA'
$ws.Cells.Item(81, 5).Value = 'NA'
$ws.Cells.Item(81, 6).Value = 'NA'
$ws.Cells.Item(81, 7).Value = 'F, W, SP'
$ws.Cells.Item(82, 3).Value = 'Admission to the MS Program in Psychology; PSY 535; PSY 555; PSY 556; and PSY 560.'
$ws.Cells.Item(82, 4).Value = 'NA'
$ws.Cells.Item(82, 5).Value = 'NA'
$ws.Cells.Item(82, 6).Value = 'NA'
$ws.Cells.Item(82, 7).Value = 'SP'
$ws.Cells.Item(83, 4).Value = 'NA'
$ws.Cells.Item(83, 5).Value = 'NA'
$ws.Cells.Item(83, 6).Value = 'NA'
$ws.Cells.Item(83, 7).Value = 'F'
$ws.Cells.Item(84, 3).Value = 'Graduate standing; PSY 560; PSY 565; and PSY 574.'
$ws.Cells.Item(84, 4).Value = 'NA'
$ws.Cells.Item(84, 5).Value = 'NA'
$ws.Cells.Item(84, 6).Value = 'NA'
$ws.Cells.Item(84, 7).Value = 'W'
$ws.Cells.Item(85, 3).Value = 'PSY 585 and advancement to candidacy.'
$ws.Cells.Item(85, 4).Value = 'NA'
$ws.Cells.Item(85, 5).Value = 'NA'
$ws.Cells.Item(85, 6).Value = 'NA'
$ws.Cells.Item(85, 7).Value = 'F, W, SP'
